$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new boolean column G with header-less boolean data for rows 1-7
$ws.Range("G1").Value = $true
$ws.Range("G2").Value = $false
$ws.Range("G3").Value = $false
$ws.Range("G4").Value = $false
$ws.Range("G5").Value = $false
$ws.Range("G6").Value = $false
$ws.Range("G7").Value = $false

# Update the selection to match the target (N7)
$ws.Range("N7").Select()
